$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Zero out the "base_damage_mod_bonus" column (AI) for rows 3-60.
# (AI2 is already 0 and is left untouched, matching the source diff.)
for ($row = 3; $row -le 60; $row++) {
    $ws.Cells.Item($row, 35).Value = 0
}

# Update the sheet's active selection to reflect the new focus area
# (AI2:AI60, active cell AI2) as captured by the diff.
[void]$ws.Range("AI2:AI60").Select()
